$p = $ppt.ActivePresentation

function Update-DateShape($shapes) {
  for ($i = 1; $i -le $shapes.Count; $i++) {
    $sh = $shapes.Item($i)
    if ($sh.HasTextFrame) {
      $tr = $sh.TextFrame.TextRange
      if ($tr.Text -eq "8/4/13") {
        $tr.Text = "8/13/13"
      }
    }
  }
}

# Slide master date placeholder
$m = $p.SlideMaster
Update-DateShape $m.Shapes

# Every slide layout's date placeholder
for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
  $layout = $m.CustomLayouts.Item($li)
  Update-DateShape $layout.Shapes
}
